$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 previously was an empty filler row (no entry yet). Add the new
# time-tracking entry: 7 March 2016, 2 hours, "Laatste fixes voor
# nieuws/stijl/menu". Pull the cell formatting from row 14 (the last
# populated data row) so the new cells match the existing date/number/text
# styling used throughout the table.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A15").Value2 = 42436
$ws.Range("B15").Value2 = 2
$ws.Range("C15").Value = "Laatste fixes voor nieuws/stijl/menu"

# The active selection moves to C16 after the entry is made.
$ws.Range("C16").Select()
